$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the memory label in A2 ("Memory4" -> "Memory2")
$ws.Range("A2").Value = "Memory2"

# Update the book id value in D2
$ws.Range("D2").Value = 178367343

# Move the active selection from A2 to D2
$ws.Range("D2").Select()
